$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (B: 8.59 -> 7.49, D: 21.79 -> 19.59, E: 10.79 -> 8.59)
# Note: the ColumnWidth COM property snaps to whole-pixel increments (1/7 character
# units at the default Calibri-11 max-digit-width of 7px), so these inputs are chosen
# to land the persisted <col width="..."> as close as possible to the target values.
$ws.Columns.Item(2).ColumnWidth = 6.714285714285714
$ws.Columns.Item(4).ColumnWidth = 18.857142857142854
$ws.Columns.Item(5).ColumnWidth = 7.857142857142858

$data = @(
    @{ Row = 2;  B = "Borg";     C = "Hague";      D = "bhague0@360.cn";              E = "Male";   F = "D" },
    @{ Row = 3;  B = "Chuck";    C = "Drover";     D = "cdrover1@goodreads.com";       E = "Male";   F = "A" },
    @{ Row = 4;  B = "Wade";     C = "Chipchase";  D = "wchipchase2@discovery.com";    E = "Male";   F = "B" },
    @{ Row = 5;  B = "Arlen";    C = "Praill";     D = "apraill3@apache.org";          E = "Female"; F = "E" },
    @{ Row = 6;  B = "Harper";   C = "Hanaford";   D = "hhanaford4@youtube.com";       E = "Male";   F = "F" },
    @{ Row = 7;  B = "Lawton";   C = "Laskey";     D = "llaskey5@nps.gov";             E = "Male";   F = "A" },
    @{ Row = 8;  B = "Hattie";   C = "Chessil";    D = "hchessil6@pinterest.com";      E = "Female"; F = "E" },
    @{ Row = 9;  B = "Paddy";    C = "Petters";    D = "ppetters7@spotify.com";        E = "Male";   F = "D" },
    @{ Row = 10; B = "Karlotte"; C = "Thomton";    D = "kthomton8@miibeian.gov.cn";    E = "Female"; F = "E" },
    @{ Row = 11; B = "Nial";     C = "Thurner";    D = "nthurner9@theguardian.com";    E = "Male";   F = "A" }
)

foreach ($rec in $data) {
    $r = $rec.Row
    $ws.Cells.Item($r, 2).Value = $rec.B
    $ws.Cells.Item($r, 3).Value = $rec.C
    $ws.Cells.Item($r, 4).Value = $rec.D
    $ws.Cells.Item($r, 5).Value = $rec.E
    $ws.Cells.Item($r, 6).Value = $rec.F
}
